$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.833.87'
$ws.Range('E2').Value = '  +4.33%  '
$ws.Range('D3').Value = '2.667.46'
$ws.Range('E3').Value = '  +7.68%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '113.44'
$ws.Range('E5').Value = '  +9.44%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '325.43'
$ws.Range('E6').Value = '  +2.98%  '
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.553'
$ws.Range('E9').Value = '  +3.99%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.81'
$ws.Range('E10').Value = '  +5.93%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.13'
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0822'
$ws.Range('E12').Value = '  +3.33%  '
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.35'
$ws.Range('E14').Value = '  +5.12%  '
$ws.Range('D15').Value = '3.087.54'
$ws.Range('E15').Value = '  +7.74%  '
$ws.Range('D16').Value = '2.661.19'
$ws.Range('E16').Value = '  +6.33%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.873'
$ws.Range('E17').Value = '  +6.54%  '
$ws.Range('D18').Value = '49.755.29'
$ws.Range('E18').Value = '  +4.30%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.14'
$ws.Range('E19').Value = '  +4.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.78'
$ws.Range('E20').Value = '  +4.45%  '
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  +3.73%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '278.13'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '71.77'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('E25').Value = '  +3.75%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.80'
$ws.Range('E26').Value = '  +4.74%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  +6.09%  '
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.19'
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('E32').Value = '  +2.23%  '
$ws.Range('E33').Value = '  +5.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '19.48'
$ws.Range('E34').Value = '  +3.56%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0810'
$ws.Range('E35').Value = '  +5.91%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.05'
$ws.Range('E37').Value = '  +13.42%  '
$ws.Range('E38').Value = '  +7.88%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.16'
$ws.Range('E39').Value = '  +11.43%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '125.27'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  +2.64%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '22.60'
$ws.Range('E42').Value = '  +4.22%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.21'
$ws.Range('E43').Value = '  +0.87%  '
$ws.Range('E44').Value = '  +6.77%  '
$ws.Range('D45').Value = '2.106.03'
$ws.Range('E45').Value = '  +5.96%  '
$ws.Range('E46').Value = '  +6.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.26'
$ws.Range('E47').Value = '  +15.50%  '
$ws.Range('E48').Value = '  +8.55%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.03'
$ws.Range('E49').Value = '  +1.97%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.35'
$ws.Range('E50').Value = '  +5.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '59.43'
$ws.Range('E51').Value = '  +7.48%  '
